$p = $ppt.ActivePresentation
$cr = [char]13

# ---------------------------------------------------------------------------
# Slide 10 - "Rectangle 6": "Sender TTL" -> "session-sender TTL"
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$shp = $s10.Shapes.Item(5)
$tr = $shp.TextFrame.TextRange
$run = $tr.Paragraphs(32, 1).Runs(1, 1)
$run.Text = "|  session-sender TTL   |      MBZ                                      |"

# ---------------------------------------------------------------------------
# Slide 11 - "Content Placeholder 6"
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$shp = $s11.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange

$tr.Paragraphs(2, 1).Runs(1, 1).Text = "Counter at fixed location - offset (no TLV, Eth 18, IPv6 40, UDP 8, Seq 4, Total = 70 Byte)"
$tr.Paragraphs(3, 1).Runs(1, 1).Text = "With LM TLV " + [char]8211 + " may not be at fixed location, also deeper into the test packet at offset (Eth 18, IPv6 40, UDP 8, STAMP 44, TLV Type 4, Total = 114 Byte)"
$tr.Paragraphs(4, 1).Runs(1, 1).Text = "Also need to include other "

# Insert a new paragraph after paragraph 5, then update paragraph 5's text.
$null = $tr.Paragraphs(5, 1).InsertAfter($cr + "Hardware also not capable to write both TS and Counter in the same packet")
$tr.Paragraphs(5, 1).Runs(1, 1).Text = "Need to load the packet in write-able memory which is limited"

$tr.Paragraphs(7, 1).Runs(1, 1).Text = "Hardware also not capable to recompute UDP checksum"
$tr.Paragraphs(9, 1).Runs(1, 1).Text = "Some test packets received from one session-sender with base test packet and some with LM TLV, hence need to parse the received test packet to check if it is for delay or direct-mode loss before punting the packet"
$tr.Paragraphs(10, 1).Runs(1, 1).Text = "Hardware need to punt with receive TS or receive Counter"
$tr.Paragraphs(11, 1).Runs(1, 1).Text = "Hardware also not capable to do both for the same packet"
$tr.Paragraphs(12, 1).Runs(1, 1).Text = "Separate UDP port + LM message format eliminate the complexity in Hardware"

# ---------------------------------------------------------------------------
# Slide 12 - "Text Placeholder 1": resize box + demote/shrink RFC bullets
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$shp = $s12.Shapes.Item(16)
$shp.Height = 919804 / 12700
$tr = $shp.TextFrame.TextRange
for ($i = 2; $i -le 4; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.IndentLevel = 2
    $para.Font.Size = 10
}

# ---------------------------------------------------------------------------
# Slide 4 - "Content Placeholder 2": "Extensions not specific to SR?" -> "Extensions specific to SR?"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp = $s4.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Paragraphs(4, 1).Runs(1, 1).Text = "Extensions specific to SR?"

# ---------------------------------------------------------------------------
# Slide 6 - "Content Placeholder 2"
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shp = $s6.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$tr.Paragraphs(1, 1).Runs(1, 1).Text = "Two-way measurement mode"

$para2 = $tr.Paragraphs(2, 1)
$para2.IndentLevel = 2
$para2.Runs(1, 1).Text = "Reflector needs to send reply on the same link (symmetric delay on forward and reverse link)"

$tr.Paragraphs(3, 1).Runs(1, 1).Text = "No way of knowing if one-way or two-way mode from the STAMP test packet"
$tr.Paragraphs(4, 1).Runs(1, 1).Text = "Not scalable to configure for each (session id, source-address) on session-reflector (can have an order of 1K links)"

# Remove the old "Reflector node may have PTP clock sync..." bullet entirely.
$tr.Paragraphs(5, 1).Delete()

$tr.Paragraphs(5, 1).Runs(1, 1).Text = "Cannot always send reply on the same incoming interface as the STAMP test packet reply may need to be IP routed"

# ---------------------------------------------------------------------------
# Slide 8 - "Content Placeholder 2": reposition/resize + shrink font + reword
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$shp = $s8.Shapes.Item(2)
$shp.Left = 533400 / 12700
$shp.Top = 1047750 / 12700
$shp.Width = 8077200 / 12700

$tr = $shp.TextFrame.TextRange
for ($i = 1; $i -le 6; $i++) {
    $tr.Paragraphs($i, 1).Font.Size = 16
}

$tr.Paragraphs(1, 1).Runs(3, 1).Text = " SR Policy, reply test packet needs to be sent on the reverse SR Policy"
$tr.Paragraphs(4, 1).Runs(1, 1).Text = "No signaling in SR, possible to use PCE"
$tr.Paragraphs(5, 1).Runs(1, 1).Text = "Need per session state on session-reflector node to store reverse paths (each session-id, source-address) " + [char]8211 + " order of 10Ks SR Policy (that can have active and standby candidate-paths and each can have multiple segment-lists)"
